# Update workbook text to reflect the new release version:
#   "mines - January 30 (built on February 02 2026 12.49.33 EST)"
# ->
#   "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# --- Sheet "About" ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: " + $newVersion

$wsAbout.Range("A6").Value = "Recommended Citation:  " + '"' + "Global Energy Monitor, Coal mine boundaries and methane sources for Changcun Coal Mine, China, M3707, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- Sheet "Boundaries and methane sources" ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$usedRows = $wsData.UsedRange.Rows.Count
for ($r = 2; $r -le $usedRows; $r++) {
    $cell = $wsData.Cells.Item($r, 19)  # column S = build_version
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value2 = $newVersion
    }
}
